# Fruta / hortaliza, semanal
# Weekly refresh of "Vega Monumental Concepción - Mandarina" data:
# two brand-new rows are inserted at the top of the existing block
# (rows 267-268), pushing the previously existing rows 267-283 down
# by two (to rows 269-285).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant columns (A..C, E..J) shared by every record in this block.
$colA = 11
$colB = "Vega Monumental Concepción"
$colC = "Bíobío"
$colE = 8
$colF = "Fruta"
$colG = 100102
$colH = "Cítricos"
$colI = 100102004
$colJ = "Mandarina"

# Target state for rows 267..285, keyed by row number:
# D (fecha-serial), K (variedad), L (calidad), M (volumen),
# N (precio minimo), O (precio maximo), P (precio promedio ponderado),
# Q (unidad comercializacion), R (origen), S (precio $/Kg), T (kg/unidad)
$rows = @{
    267 = @(45223, "Murcott",     "Primera", 100, 8000,  8000,  8000,  "`$/bandeja 18 kilos", "Región de O'Higgins", 444,  18)
    268 = @(45223, "Murcott",     "Segunda", 100, 7000,  7000,  7000,  "`$/bandeja 18 kilos", "Región de O'Higgins", 389,  18)
    269 = @(44468, "Murcott",     "Primera", 200, 7000,  7500,  7250,  "`$/bandeja 18 kilos", "Provincia de Limarí", 403,  18)
    270 = @(44468, "Murcott",     "Segunda", 100, 6500,  6500,  6500,  "`$/bandeja 18 kilos", "Provincia de Limarí", 361,  18)
    271 = @(44778, "Clemenuless", "Primera", 450, 7500,  8000,  7722,  "`$/caja 18 kilos",    "Región de O'Higgins", 429,  18)
    272 = @(44778, "Clemenuless", "Segunda", 300, 6500,  6500,  6500,  "`$/caja 18 kilos",    "Región de O'Higgins", 361,  18)
    273 = @(44754, "Clementina",  "Primera", 250, 8000,  8500,  8300,  "`$/caja 18 kilos",    "Región de O'Higgins", 461,  18)
    274 = @(44754, "Clementina",  "Segunda", 220, 6500,  7000,  6727,  "`$/caja 18 kilos",    "Región de O'Higgins", 374,  18)
    275 = @(44874, "Murcott",     "Primera", 200, 8000,  9000,  8500,  "`$/bandeja 18 kilos", "Región de O'Higgins", 472,  18)
    276 = @(44775, "Clemenuless", "Primera", 100, 8000,  9000,  8500,  "`$/bandeja 18 kilos", "Región de O'Higgins", 472,  18)
    277 = @(44775, "Clemenuless", "Segunda", 50,  7000,  7000,  7000,  "`$/bandeja 18 kilos", "Región de O'Higgins", 389,  18)
    278 = @(45063, "Murcott",     "Primera", 220, 10000, 11000, 10455, "`$/bandeja 10 kilos", "Provincia de Limarí", 1046, 10)
    279 = @(44365, "Clementina",  "Primera", 200, 8000,  9000,  8500,  "`$/bandeja 10 kilos", "Provincia de Limarí", 850,  10)
    280 = @(44365, "Clementina",  "Segunda", 100, 7000,  7000,  7000,  "`$/bandeja 10 kilos", "Provincia de Limarí", 700,  10)
    281 = @(44818, "Murcott",     "Primera", 100, 7500,  8000,  7750,  "`$/bandeja 18 kilos", "Región de O'Higgins", 431,  18)
    282 = @(44818, "Murcott",     "Segunda", 50,  6500,  6500,  6500,  "`$/bandeja 18 kilos", "Región de O'Higgins", 361,  18)
    283 = @(44341, "Clemenuless", "Primera", 200, 10000, 11000, 10500, "`$/bandeja 10 kilos", "Provincia de Limarí", 1050, 10)
    284 = @(45216, "Murcott",     "Primera", 250, 8000,  9000,  8400,  "`$/bandeja 18 kilos", "Región de O'Higgins", 467,  18)
    285 = @(44392, "Clemenuless", "Primera", 200, 6500,  7000,  6750,  "`$/bandeja 10 kilos", "Provincia de Limarí", 675,  10)
}

foreach ($r in 267..285) {
    $vals = $rows[$r]

    # Columns A-C and E-J are identical for every row in this block; rows
    # 267-283 already carry them, but write them unconditionally so the two
    # brand-new rows (284, 285) get them too.
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ

    # Column D (Fecha) needs the date/time number format applied, like the
    # rest of the column.
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 11).Value = $vals[1]
    $ws.Cells.Item($r, 12).Value = $vals[2]
    $ws.Cells.Item($r, 13).Value = $vals[3]
    $ws.Cells.Item($r, 14).Value = $vals[4]
    $ws.Cells.Item($r, 15).Value = $vals[5]
    $ws.Cells.Item($r, 16).Value = $vals[6]
    $ws.Cells.Item($r, 17).Value = $vals[7]
    $ws.Cells.Item($r, 18).Value = $vals[8]
    $ws.Cells.Item($r, 19).Value = $vals[9]
    $ws.Cells.Item($r, 20).Value = $vals[10]
}
